$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit shuffles the species-record rows 14-22 (each row keeps its
# full set of column values together, only columns A,B,D,E,F,G,H,Q,R
# differ between the rows - the rest of the columns are identical for
# every one of these rows). Below are the new values for each row after
# the shuffle.

$rows = @{
    14 = @{ A = 111251423; B = 77677;  D = 'VU'; E = 1249; F = 'Norsk näverlav';     G = 'Platismatia norvegica';   H = '(Lynge) W.L.Culb. & C.F.Culb.'; Q = 460188.7895233887; R = 7164860.82616597 }
    15 = @{ A = 111251432; B = 81248;  D = 'NT'; E = 1312; F = 'Gammelgransskål';    G = 'Pseudographis pinicola';  H = '(Nyl.) Rehm';                   Q = 460622.5513675315; R = 7165027.330594921 }
    16 = @{ A = 111251420; B = 77677;  D = 'VU'; E = 1249; F = 'Norsk näverlav';     G = 'Platismatia norvegica';   H = '(Lynge) W.L.Culb. & C.F.Culb.'; Q = 460243.4530616797; R = 7164800.429238674 }
    17 = @{ A = 111251407; B = 73696;  D = 'NT'; E = 6440; F = 'Vitgrynig nållav';   G = 'Chaenotheca subroscida';  H = '(Eitner) Zahlbr.';              Q = 460240.5118381025; R = 7164805.620072429 }
    18 = @{ A = 111251428; B = 89423;  D = 'NT'; E = 5432; F = 'Granticka';          G = 'Porodaedalea chrysoloma'; H = '(Fr.) Fiasson & Niemelä';       Q = 460445.0942901828; R = 7164835.148113105 }
    19 = @{ A = 111251402; B = 73696;  D = 'NT'; E = 6440; F = 'Vitgrynig nållav';   G = 'Chaenotheca subroscida';  H = '(Eitner) Zahlbr.';              Q = 460212.3128264685; R = 7164818.870384302 }
    20 = @{ A = 111251434; B = 78612;  D = 'LC'; E = 6464; F = 'Luddlav';            G = 'Nephroma resupinatum';    H = '(L.) Ach.';                     Q = 460452.9763639791; R = 7164846.208533676 }
    21 = @{ A = 111251437; B = 78611;  D = 'LC'; E = 6463; F = 'Bårdlav';            G = 'Nephroma parile';         H = '(Ach.) Ach.';                   Q = 460452.9763639791; R = 7164846.208533676 }
    22 = @{ A = 111251430; B = 77515;  D = 'NT'; E = 6425; F = 'Garnlav';            G = 'Alectoria sarmentosa';    H = '(Ach.) Ach.';                   Q = 460188.8289468794; R = 7164863.831099218 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("A$r").Value = $vals.A
    $ws.Range("B$r").Value = $vals.B
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("E$r").Value = $vals.E
    $ws.Range("F$r").Value = $vals.F
    $ws.Range("G$r").Value = $vals.G
    $ws.Range("H$r").Value = $vals.H
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
}
